$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Neo4j query text stored in cell B2 (row 2) of the "startup" sheet.
# Single-quoted here-string keeps backticks and special characters literal.
$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['BAI']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Range("B2").Value = $newQuery

# The longer query text requires a taller row to display, matching the
# updated layout in the workbook (row height 201.5 -> 279).
$ws.Rows.Item(2).RowHeight = 279
